$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Range("B11")

# Cell B11 held the shared string "R40"; it now holds the text string "1",
# kept in the same (General) style/format it already had.
#
# A plain `.Value = "1"` gets auto-coerced by Excel to the *number* 1
# because the cell's format is General - so the literal text "1" needs to
# be produced some other way, without touching B11's own NumberFormat
# (which would otherwise permanently grow the workbook's style table with
# an unused Text xf, since styles are never garbage collected).
#
# TEXT() always returns a string, so evaluating it on a scratch cell and
# pasting just the computed value back onto B11 deposits the literal text
# "1" there while leaving B11's style/format completely untouched.
$scratch = $ws.Range("ZZ9000")
$scratch.Formula = '=TEXT(1,"0")'
$scratch.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()
$excel.CutCopyMode = $false
